# Update the second Yahboom GitHub URL on the "Offline programming" slide:
#   https://github.com/YahboomTechnology/YB_IR.
# becomes
#   https://github.com/YahboomTechnology/Yahboom_IR
# (the trailing "." that used to live in its own run is absorbed, so that
# run disappears and the following sentence still starts with ". Then ...")

$p = $ppt.ActivePresentation

$oldUrl = "https://github.com/YahboomTechnology/YB_IR."
$newUrl = "https://github.com/YahboomTechnology/Yahboom_IR"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            $full = $tr.Text
            $idx = $full.IndexOf($oldUrl)
            while ($idx -ge 0) {
                $sub = $tr.Characters($idx + 1, $oldUrl.Length)
                $sub.Text = $newUrl
                $full = $tr.Text
                $idx = $full.IndexOf($oldUrl)
            }
        }
    }
}
